$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: shift the "Tecnologia dos Materiais - MEC-1A" entry from column C
# to columns E and F, leaving C and D as "-"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "Tecnologia dos Materiais - MEC-1A"
$ws.Range("F4").Value = "Tecnologia dos Materiais - MEC-1A"

# Row 6: clear the "Tecnologia dos Materiais - MEC-1A" entry from column C
$ws.Range("C6").Value = "-"
